$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (rows 214-219), matching columns A..K:
# A: FECHA_OPERACION (date serial), B: CONTRATO, C: CONTRATO_RECTIFICACION,
# D: CONTRATO_ANULACION, E: CONTRATO_PRECIO_HECHO, F: FIJACION,
# G: FIJACION_RECTIFICACION, H: FIJACION_ANULACION, I: FIJACIONES,
# J: TOTAL, K: PRODUCTO
$rows = @(
    @(45950, 25070,     40,   0, 25110,     4340, 0, 0, 4340,     29450,    "CEBADA"),
    @(45951, 32288.12,   0,   400, 31888.12, 500,  0, 0, 500,     32388.12, "CEBADA"),
    @(45952, 23990.26, 4240,  0, 28230.26,  1500, 0, 0, 1500,    29730.26, "CEBADA"),
    @(45953, 10986,     500,  0, 11486,     30,   0, 0, 30,      11516,    "CEBADA"),
    @(45954, 14826.67,  710,  0, 15536.67,  51.18000000000001, 0, 0, 51.18000000000001, 15587.85, "CEBADA"),
    @(45957, 450,        0,   0, 450,       7500, 0, 0, 7500,    7950,     "CEBADA")
)

$startRow = 214
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd"

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
    $ws.Cells.Item($r, 10).Value = $data[9]
    $ws.Cells.Item($r, 11).Value = $data[10]
}
